# Adds the three newly-scraped practice/service lookup sheets -- "kubrick",
# "iqvia" and "cambridge_consultants" -- to the end of the workbook. Each one
# gets a "practices"/"services" header (bold, bordered, center/top-aligned)
# followed by the scraped practice -> service rows, matching the existing
# scrape sheets (cognizant / scrape_bettergov / scrape_capco / infosys).
#
# New sheets are created by copying an existing scrape sheet ("infosys") and
# then overwriting its cells in place: this clones the sheet-level template
# (sheetPr/outlinePr, pageMargins, sheetFormatPr, and the headers cellXf)
# exactly, rather than relying on a freshly-Add()-ed sheets different
# defaults, and leaves the donor sheet itself completely untouched.

$wb = $excel.ActiveWorkbook
$donor = $wb.Worksheets.Item("infosys")
$donorRows = $donor.UsedRange.Rows.Count

# ---- sheet "kubrick" (A1:B44) ----
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$donor.Copy($null, $lastSheet)
$ws = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws.Name = "kubrick"

$kubrickData = @(
    ("practices", "services"),
    ("Data Product Management", "Business Analysis & Business Process Mining"),
    ("Data Product Management", "Requirements Engineering"),
    ("Data Product Management", "Agile Product Management"),
    ("Data Product Management", "Change & Risk Management"),
    ("Data Product Management", "Data Product Testing & Validation"),
    ("Data Product Management", "Data Product Launch & Deployment"),
    ("Data Product Management", "Data Product Lifecycle Management"),
    ("Data Engineering", "Modern Data Architecture & Data Modelling"),
    ("Data Engineering", "Data Connectivity and Integration"),
    ("Data Engineering", "Cloud Data Warehouse and Lake Development"),
    ("Data Engineering", "DataOps"),
    ("Data Engineering", "Data Orchestration"),
    ("Data Engineering", "Data Streaming" + [char]0x200B + ""),
    ("Data Engineering", "Scalability & Performance Optimisation"),
    ("Data Engineering", "Database Design & Development" + [char]0x200B + ""),
    ("Data & AI Governance", "Data Governance Frameworks & Policy"),
    ("Data & AI Governance", "Policy-as-Code"),
    ("Data & AI Governance", "Master Data & Reference Data Management"),
    ("Data & AI Governance", "Data Quality Management"),
    ("Data & AI Governance", "Data Cataloging & Lineage"),
    ("Data & AI Governance", "Data Privacy & Compliance"),
    ("Data & AI Governance", "Data Domain Modelling"),
    ("Advanced Analytics", "Data Analysis & Insights" + [char]0x200B + ""),
    ("Advanced Analytics", "Advanced Data Visualisation"),
    ("Advanced Analytics", "Decision Intelligence" + [char]0x200B + ""),
    ("Advanced Analytics", "Knowledge Graph Development" + [char]0x200B + ""),
    ("Advanced Analytics", "Digital Twins" + [char]0x200B + ""),
    ("Advanced Analytics", "Data Storytelling" + [char]0x200B + ""),
    ("Advanced Analytics", "Self-service Enablement"),
    ("GenAI & MLOps", "Feature Engineering" + [char]0x200B + ""),
    ("GenAI & MLOps", "Model Development" + [char]0x200B + ""),
    ("GenAI & MLOps", "ML Engineering " + [char]0x200B + ""),
    ("GenAI & MLOps", "LLMOps" + [char]0x200B + ""),
    ("GenAI & MLOps", "LLM Integration & Fine Tuning" + [char]0x200B + ""),
    ("GenAI & MLOps", "Prompt Engineering" + [char]0x200B + ""),
    ("GenAI & MLOps", "AI Ethics & Compliance"),
    ("Cloud", "Cloud Design & Deployment" + [char]0x200B + ""),
    ("Cloud", "Cloud Migration" + [char]0x200B + ""),
    ("Cloud", "Cloud Infra Optimization" + [char]0x200B + ""),
    ("Cloud", "CI/CD" + [char]0x200B + ""),
    ("Cloud", "FinOps & Sustainability" + [char]0x200B + ""),
    ("Cloud", "Cloud Security & Compliance" + [char]0x200B + ""),
    ("Cloud", "SRE")
)

for ($i = 0; $i -lt $kubrickData.Length; $i++) {
    $r = $i + 1
    $pair = $kubrickData[$i]
    $ws.Cells.Item($r, 1).Value = $pair[0]
    $ws.Cells.Item($r, 2).Value = $pair[1]
}

# The donor sheet had more rows than this dataset needs -- clear the leftover
# copied cells (content + the header-only style never reaches past row 1, so
# this is a plain content clear) so the sheet dimension shrinks to match.
$kubrickRows = $kubrickData.Length
if ($donorRows -gt $kubrickRows) {
    $extra = $ws.Range($ws.Cells.Item($kubrickRows + 1, 1), $ws.Cells.Item($donorRows, 2))
    $extra.Clear()
}

# ---- sheet "iqvia" (A1:B14) ----
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$donor.Copy($null, $lastSheet)
$ws = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws.Name = "iqvia"

$iqviaData = @(
    ("practices", "services"),
    ("Clinical Research", "Monitoring"),
    ("Clinical Research", "ClinicalOperations"),
    ("Clinical Research", "Clinical Project Management & Leadership"),
    ("Clinical Research", "Clinical DataManagement"),
    ("Clinical Research", "StatisticalServices"),
    ("Technology & Analytics", "InformationSecurity"),
    ("Technology & Analytics", "DevOps"),
    ("Technology & Analytics", "Software Development Engineering"),
    ("Technology & Analytics", "Software QA and Test Engineering"),
    ("Technology & Analytics", "AI and Machine Learning"),
    ("Consulting", "Commercial Consulting"),
    ("Consulting", "Real World Consulting"),
    ("Consulting", "Technical Consulting")
)

for ($i = 0; $i -lt $iqviaData.Length; $i++) {
    $r = $i + 1
    $pair = $iqviaData[$i]
    $ws.Cells.Item($r, 1).Value = $pair[0]
    $ws.Cells.Item($r, 2).Value = $pair[1]
}

# The donor sheet had more rows than this dataset needs -- clear the leftover
# copied cells (content + the header-only style never reaches past row 1, so
# this is a plain content clear) so the sheet dimension shrinks to match.
$iqviaRows = $iqviaData.Length
if ($donorRows -gt $iqviaRows) {
    $extra = $ws.Range($ws.Cells.Item($iqviaRows + 1, 1), $ws.Cells.Item($donorRows, 2))
    $extra.Clear()
}

# ---- sheet "cambridge_consultants" (A1:B31) ----
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$donor.Copy($null, $lastSheet)
$ws = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws.Name = "cambridge_consultants"

$cambridge_consultantsData = @(
    ("practices", "services"),
    ("5G and wireless connectivity", "Radio systems"),
    ("5G and wireless connectivity", "Connectivity and IOT"),
    ("5G and wireless connectivity", "Digital signal processing"),
    ("5G and wireless connectivity", "AI and analytics"),
    ("5G and wireless connectivity", "Strategic advice"),
    ("Advanced computing", "ASICs and electronics"),
    ("Advanced computing", "Optics and photonics"),
    ("Advanced computing", "Physical sciences"),
    ("Advanced computing", "AI and analytics"),
    ("AI and data analytics", "AI and analytics"),
    ("AI and data analytics", "Simulation"),
    ("AI and data analytics", "Connectivity and IoT"),
    ("AI and data analytics", "Sensing"),
    ("AI and data analytics", "Electronics and ASICs"),
    ("Biotechnology", "Cell and gene therapy"),
    ("Biotechnology", "Synthetic biology"),
    ("Biotechnology", "AI and analytics"),
    ("Biotechnology", "Physical sciences"),
    ("Biotechnology", "Simulation"),
    ("Biotechnology", "Strategic advice"),
    ("Digital transformation", "Digital services"),
    ("Digital transformation", "Digital security"),
    ("Digital transformation", "AI and analytics"),
    ("Digital transformation", "Extended reality (XR)"),
    ("Digital transformation", "User experience (UX)"),
    ("Digital transformation", "Connectivity and IOT"),
    ("Quantum technology", "Optics and photonics"),
    ("Quantum technology", "Physical sciences"),
    ("Quantum technology", "Product realisation"),
    ("Quantum technology", "Strategic advice")
)

for ($i = 0; $i -lt $cambridge_consultantsData.Length; $i++) {
    $r = $i + 1
    $pair = $cambridge_consultantsData[$i]
    $ws.Cells.Item($r, 1).Value = $pair[0]
    $ws.Cells.Item($r, 2).Value = $pair[1]
}

# The donor sheet had more rows than this dataset needs -- clear the leftover
# copied cells (content + the header-only style never reaches past row 1, so
# this is a plain content clear) so the sheet dimension shrinks to match.
$cambridge_consultantsRows = $cambridge_consultantsData.Length
if ($donorRows -gt $cambridge_consultantsRows) {
    $extra = $ws.Range($ws.Cells.Item($cambridge_consultantsRows + 1, 1), $ws.Cells.Item($donorRows, 2))
    $extra.Clear()
}

# Restore the original active sheet/selection (Sheet1) so workbook-level view
# state is unchanged, matching the diff (which touches only the sheet list).
[void]$wb.Worksheets.Item("Sheet1").Activate()
[void]$wb.Worksheets.Item("Sheet1").Range("A1").Select()
